# Updates the cryptos price/volume table with freshly scraped values.
# Note: several Price-column strings look numeric (e.g. "1.00", "15.40",
# "0.0000267") and must stay plain text (as in the source workbook), so a
# leading apostrophe is used to force Excel to keep them as text instead of
# normalising them into numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.205.51"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").Value = "3.178.32"
$ws.Range("E3").Value = "  -1.02%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'611.26"
$ws.Range("E5").Value = "  +1.36%  "
$ws.Range("E6").Value = "  +0.70%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "3.176.46"
$ws.Range("E8").Value = "  -1.08%  "
$ws.Range("D9").Value = "'0.545"
$ws.Range("E9").Value = "  +2.29%  "
$ws.Range("E10").Value = "  -0.96%  "
$ws.Range("D11").Value = "'5.68"
$ws.Range("E11").Value = "  -7.24%  "
$ws.Range("D12").Value = "'0.516"
$ws.Range("E12").Value = "  +0.96%  "
$ws.Range("D13").Value = "'0.0000267"
$ws.Range("E13").Value = "  -1.45%  "
$ws.Range("D14").Value = "'38.42"
$ws.Range("E14").Value = "  -2.96%  "
$ws.Range("D15").Value = "3.699.24"
$ws.Range("E15").Value = "  -1.07%  "
$ws.Range("D16").Value = "66.235.44"
$ws.Range("E16").Value = "  +0.25%  "
$ws.Range("D17").Value = "'7.41"
$ws.Range("E17").Value = "  -1.40%  "
$ws.Range("D18").Value = "3.178.32"
$ws.Range("E18").Value = "  -1.25%  "
$ws.Range("E19").Value = "  +0.90%  "
$ws.Range("D20").Value = "'510.16"
$ws.Range("E20").Value = "  -0.16%  "
$ws.Range("D21").Value = "'15.40"
$ws.Range("E21").Value = "  -0.51%  "
$ws.Range("E22").Value = "  -1.47%  "
$ws.Range("E23").Value = "  -1.24%  "
$ws.Range("D24").Value = "'14.84"
$ws.Range("E24").Value = "  -4.69%  "
$ws.Range("D25").Value = "'84.59"
$ws.Range("E25").Value = "  -0.55%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("E27").Value = "  -0.22%  "
$ws.Range("D28").Value = "'9.14"
$ws.Range("E28").Value = "  -2.61%  "
$ws.Range("D29").Value = "'2.39"
$ws.Range("E29").Value = "  +4.39%  "
$ws.Range("D30").Value = "'3.01"
$ws.Range("E30").Value = "  +4.57%  "
$ws.Range("D31").Value = "'7.17"
$ws.Range("E31").Value = "  +4.78%  "
$ws.Range("D32").Value = "'27.98"
$ws.Range("E32").Value = "  -0.68%  "
$ws.Range("E33").Value = "  +0.10%  "
$ws.Range("E34").Value = "  -1.88%  "
$ws.Range("D35").Value = "'6.51"
$ws.Range("E35").Value = "  -1.30%  "
$ws.Range("D36").Value = "'506.36"
$ws.Range("E36").Value = "  +3.89%  "
$ws.Range("D37").Value = "'55.02"
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("E38").Value = "  -2.85%  "
$ws.Range("E39").Value = "  +0.10%  "
$ws.Range("E40").Value = "  +5.03%  "
$ws.Range("E41").Value = "  -1.55%  "
$ws.Range("D42").Value = "0.0₃0685"
$ws.Range("E42").Value = "  +6.49%  "
$ws.Range("D43").Value = "'2.85"
$ws.Range("E43").Value = "  -3.90%  "
$ws.Range("E44").Value = "  -2.44%  "
$ws.Range("D45").Value = "'2.45"
$ws.Range("E45").Value = "  -0.58%  "
$ws.Range("D46").Value = "2.828.18"
$ws.Range("E46").Value = "  -4.33%  "
$ws.Range("E47").Value = "  -2.11%  "
$ws.Range("B48").Value = "ThetaToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D48").Value = "'2.37"
$ws.Range("E48").Value = "  +2.25%  "
$ws.Range("B49").Value = "USDe"
$ws.Range("C49").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D49").Value = "'0.999"
$ws.Range("E49").Value = "  -0.11%  "
$ws.Range("E50").Value = "  +0.32%  "
$ws.Range("E51").Value = "  +7.21%  "
